# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 2
    $ws.Range("F3").Value = 297
    $ws.Range("F4").Value = 207
    $ws.Range("F5").Value = 2514
    $ws.Range("F6").Value = 1819
    $ws.Range("F7").Value = 348

    if ($sheetName -eq "展览") {
        $ws.Range("F8").Value = 103
        $ws.Range("F9").Value = 872
    } else {
        $ws.Range("F9").Value = 103
        $ws.Range("F10").Value = 872
    }
}
